# Add new Beatrix Potter books to the BooksLibrary sheet
# (fix typo on the existing "Peter Rabbit" row and insert 7 new rows
#  of Beatrix Potter titles just below it, highlighted in red).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 7 new Beatrix Potter rows right after row 6
# ("The Tale of Peter Rabbit" / row 6), pushing the rest of the
# table (old rows 7-24) down to rows 14-31.
$ws.Rows("7:13").Insert()

# New row 7
$ws.Range("A7").Value = "The Tale of Two Bad Mice"
$ws.Range("B7").Value = "BEATRIX POTTER"
$ws.Range("C7").Value = "https://www.gutenberg.org/ebooks/45264"

# New row 8
$ws.Range("A8").Value = "The Tale of the Pie and the Patty Pan"
$ws.Range("B8").Value = "BEATRIX POTTER"
$ws.Range("C8").Value = "https://www.gutenberg.org/ebooks/15234"

# New rows 9-13 - titles first
$ws.Range("A9").Value = "The Tale of Mrs. Tittlemouse"
$ws.Range("B9").Value = "BEATRIX POTTER"

$ws.Range("A10").Value = "The Tale of Mrs. Tiggy-Winkle"
$ws.Range("B10").Value = "BEATRIX POTTER"

$ws.Range("A11").Value = "The Tale of Ginger and Pickles"
$ws.Range("B11").Value = "BEATRIX POTTER"

$ws.Range("A12").Value = "The Story of Miss Moppet"
$ws.Range("B12").Value = "BEATRIX POTTER"

$ws.Range("A13").Value = "The Story of a Fierce Bad Rabbit"
$ws.Range("B13").Value = "BEATRIX POTTER"

# ... then the links for rows 9-12
$ws.Range("C9").Value = "https://www.gutenberg.org/ebooks/17089"
$ws.Range("C10").Value = "https://www.gutenberg.org/ebooks/15137"
$ws.Range("C11").Value = "https://www.gutenberg.org/ebooks/14877"
$ws.Range("C12").Value = "https://www.gutenberg.org/ebooks/14848"

# Fix the long standing typo in the existing Peter Rabbit title.
$ws.Range("A6").Value = "The Tale of Peter Rabbit"

# Row 13's link is filled in last.
$ws.Range("C13").Value = "https://www.gutenberg.org/ebooks/45265"

# Highlight the first and last new rows (7 and 13) in red, as done
# for the two new entries bracketing the batch.
$ws.Range("A7:C7").Font.Color = 255
$ws.Range("A13:C13").Font.Color = 255

# Restore a sane view: no frozen/scrolled top-left cell, and the
# last edited row (13) selected in full.
$ws.Range("A13:XFD13").Select()
